$d = $word.ActiveDocument

# 1) Replace the "m:enduserdoc" field (fldChar begin / instrText / fldChar end)
#    with a plain literal text run "{m:enduserdoc}".
#    Insert the literal text immediately before the field first (so it lands
#    in its own, unformatted run), then delete the field itself.
$f = $d.Fields.Item(1)
$fieldStart = $f.Code.Start - 1

$r = $d.Range($fieldStart, $fieldStart)
$r.InsertBefore("{m:enduserdoc}")

$f.Delete()

# 2) Prefix the error message run with the "    <---" marker.
$d.Content.Find.Execute("Invalid block: Unexpected tag m:enduserdoc at this location", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "    <---Invalid block: Unexpected tag m:enduserdoc at this location", 2)
